# Update the test data row (row 2) on the active sheet (Hoja1) of the
# collections revenue report. Column headers in row 1 stay the same;
# only the sample data in row 2 is replaced with new "test" values, and
# the numeric Amount/VAT/Total columns get populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "c4b230aff0"
$ws.Range("B2").Value = "21/06/2013"
$ws.Range("C2").Value = "Miguel Angel Canas Vaz"
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 12
